# Updated cryptos list - apply per-row Price (D) and Volume(1h) (E) updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    [PSCustomObject]@{ Row = 2; D = '25.884.93'; E = '  -1.13%  ' }
    [PSCustomObject]@{ Row = 3; D = '1.639.49'; E = '  -1.01%  ' }
    [PSCustomObject]@{ Row = 4; D = '1.002'; E = '  -0.91%  ' }
    [PSCustomObject]@{ Row = 5; D = '215.37'; E = '  -0.47%  ' }
    [PSCustomObject]@{ Row = 6; D = '0.5038'; E = '  -1.70%  ' }
    [PSCustomObject]@{ Row = 7; D = '1.003'; E = '  -0.95%  ' }
    [PSCustomObject]@{ Row = 8; D = '0.2573'; E = '  -0.94%  ' }
    [PSCustomObject]@{ Row = 9; D = '0.06386'; E = '  -0.93%  ' }
    [PSCustomObject]@{ Row = 10; D = '19.51'; E = '  -1.31%  ' }
    [PSCustomObject]@{ Row = 11; D = '0.07739'; E = '  -0.85%  ' }
    [PSCustomObject]@{ Row = 12; D = '1.645.51'; E = '  -0.56%  ' }
    [PSCustomObject]@{ Row = 13; D = '4.258'; E = '  -0.63%  ' }
    [PSCustomObject]@{ Row = 14; D = '1.866.03'; E = '  -0.91%  ' }
    [PSCustomObject]@{ Row = 15; D = '0.5460'; E = '  -0.54%  ' }
    [PSCustomObject]@{ Row = 16; D = '0.0₅7894'; E = '  -1.43%  ' }
    [PSCustomObject]@{ Row = 17; D = '64.23'; E = '  +0.32%  ' }
    [PSCustomObject]@{ Row = 18; D = '25.903.31'; E = '  -1.10%  ' }
    [PSCustomObject]@{ Row = 19; D = '1.002'; E = '  -0.84%  ' }
    [PSCustomObject]@{ Row = 20; D = '201.76'; E = '  -3.20%  ' }
    [PSCustomObject]@{ Row = 21; D = '4.384'; E = '  -0.41%  ' }
    [PSCustomObject]@{ Row = 22; D = '9.887'; E = '  -2.03%  ' }
    [PSCustomObject]@{ Row = 23; D = '5.975'; E = '  -1.10%  ' }
    [PSCustomObject]@{ Row = 24; D = '1.003'; E = '  -0.89%  ' }
    [PSCustomObject]@{ Row = 25; D = '1.867'; E = '  +0.83%  ' }
    [PSCustomObject]@{ Row = 26; D = '140.99'; E = '  -2.35%  ' }
    [PSCustomObject]@{ Row = 27; D = '0.1136'; E = '  -3.00%  ' }
    [PSCustomObject]@{ Row = 28; D = '15.64'; E = '  -1.23%  ' }
    [PSCustomObject]@{ Row = 29; D = '6.775'; E = '  -2.62%  ' }
    [PSCustomObject]@{ Row = 30; D = '1.243'; E = '  -0.07%  ' }
    [PSCustomObject]@{ Row = 31; D = '0.04951'; E = '  -2.66%  ' }
    [PSCustomObject]@{ Row = 32; D = '3.270'; E = '  -2.31%  ' }
    [PSCustomObject]@{ Row = 33; D = '3.199'; E = '  -1.37%  ' }
    [PSCustomObject]@{ Row = 34; D = '1.546'; E = '  -0.47%  ' }
    [PSCustomObject]@{ Row = 35; D = '2.363'; E = '  +0.06%  ' }
    [PSCustomObject]@{ Row = 36; D = '2.628'; E = '  -4.13%  ' }
    [PSCustomObject]@{ Row = 37; D = '0.8917'; E = '  -3.04%  ' }
    [PSCustomObject]@{ Row = 38; D = '1.147.70'; E = '  -2.23%  ' }
    [PSCustomObject]@{ Row = 39; D = '0.5592'; E = '  -2.00%  ' }
    [PSCustomObject]@{ Row = 40; D = $null; E = '  -1.27%  ' }
    [PSCustomObject]@{ Row = 41; D = '1.003'; E = '  -0.91%  ' }
    [PSCustomObject]@{ Row = 42; D = '5.688'; E = '  +0.40%  ' }
    [PSCustomObject]@{ Row = 43; D = '99.70'; E = '  -0.80%  ' }
    [PSCustomObject]@{ Row = 44; D = '0.8054'; E = '  -2.59%  ' }
    [PSCustomObject]@{ Row = 45; D = '1.777.71'; E = $null }
    [PSCustomObject]@{ Row = 46; D = '0.0₈118'; E = '  +4.83%  ' }
    [PSCustomObject]@{ Row = 47; D = '0.4526'; E = '  -0.77%  ' }
    [PSCustomObject]@{ Row = 48; D = $null; E = '  -0.23%  ' }
    [PSCustomObject]@{ Row = 49; D = '54.73'; E = $null }
    [PSCustomObject]@{ Row = 50; D = '0.05052'; E = '  -0.75%  ' }
    [PSCustomObject]@{ Row = 51; D = '1.000'; E = '  -1.04%  ' }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $cell = $ws.Cells.Item($u.Row, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
    }
    if ($null -ne $u.E) {
        $cell = $ws.Cells.Item($u.Row, 5)
        $cell.NumberFormat = "@"
        $cell.Value = $u.E
    }
}

Write-Host "Updated $($updates.Count) rows."
